$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.295.70"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "'3.248.02"
$ws.Range("E3").Value = "  +2.60%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'610.04"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").Value = "'156.90"
$ws.Range("E6").Value = "  +2.10%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'3.246.12"
$ws.Range("E8").Value = "  +2.63%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  +1.89%  "
$ws.Range("E11").Value = "  +3.02%  "
$ws.Range("E12").Value = "  -4.09%  "
$ws.Range("D14").Value = "'38.98"
$ws.Range("E14").Value = "  +1.98%  "
$ws.Range("D15").Value = "'3.785.31"
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("D16").Value = "'66.422.45"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "'3.252.51"
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("D20").Value = "'503.63"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").Value = "'15.43"
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D22").Value = "'0.752"
$ws.Range("E22").Value = "  +3.64%  "
$ws.Range("D23").Value = "'8.06"
$ws.Range("E23").Value = "  +1.01%  "
$ws.Range("D24").Value = "'14.62"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("D25").Value = "'87.22"
$ws.Range("E25").Value = "  +3.20%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'3.01"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").Value = "'9.16"
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("D30").Value = "'0.125"
$ws.Range("E30").Value = "  +42.86%  "
$ws.Range("D31").Value = "'6.98"
$ws.Range("E31").Value = "  -2.13%  "
$ws.Range("E32").Value = "  -4.22%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  -3.68%  "
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("D37").Value = "'55.46"
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("D38").Value = "'3.30"
$ws.Range("E38").Value = "  +17.90%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "'0.0₃0780"
$ws.Range("E39").Value = "  +14.87%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'493.80"
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("D41").Value = "'0.0420"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("E42").Value = "  +1.25%  "
$ws.Range("D43").Value = "'8.84"
$ws.Range("E43").Value = "  +1.29%  "
$ws.Range("E44").Value = "  +4.05%  "
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("D46").Value = "'2.981.18"
$ws.Range("E46").Value = "  +5.91%  "
$ws.Range("D47").Value = "'28.82"
$ws.Range("E47").Value = "  +3.33%  "
$ws.Range("E48").Value = "  +5.34%  "
$ws.Range("E49").Value = "  +2.56%  "
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("E51").Value = "  -2.35%  "
